$d = $word.ActiveDocument

# The document has a "default" header/footer pair and a "first page"
# header/footer pair (Sections(1).Headers/Footers Item(1) = default,
# Item(2) = first page) - both pairs carry the same two logos:
#   - the BTec logo picture in the headers
#   - the Pearson logo picture in the footers
#
# Rename the embedded picture InlineShapes:
#   BTec_Logo-Orange : image1.jpg -> image2.jpg (both headers)
#   PearsonLogo.png  : image2.png -> image1.png (both footers)

for ($i = 1; $i -le 2; $i++) {
    $hdr = $d.Sections.Item(1).Headers.Item($i)
    if ($hdr.Exists) {
        for ($s = 1; $s -le $hdr.Range.InlineShapes.Count; $s++) {
            $shp = $hdr.Range.InlineShapes.Item($s)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }

    $ftr = $d.Sections.Item(1).Footers.Item($i)
    if ($ftr.Exists) {
        for ($s = 1; $s -le $ftr.Range.InlineShapes.Count; $s++) {
            $shp = $ftr.Range.InlineShapes.Item($s)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
